# Hortaliza, Terminal La Palmera de La Serena - Ajo
# A new weekly price-report row is inserted at row 31 (pushing the
# existing rows 31..121 down to 32..122), matching the "Fruta / hortaliza,
# semanal" update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row at position 31; everything below shifts down one.
$ws.Rows("31:31").Insert()

# Populate the new row 31 with the new weekly data point.
$ws.Range("A31").Value = 8
$ws.Range("B31").Value = "Terminal La Palmera de La Serena"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44414
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 100112003
$ws.Range("G31").Value = "Ajo"
$ws.Range("H31").Value = "Chino"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 720
$ws.Range("K31").Value = 12500
$ws.Range("L31").Value = 13000
$ws.Range("M31").Value = 12750
$ws.Range("N31").Value = "$/caja 10 kilos"
$ws.Range("O31").Value = "China"
$ws.Range("P31").Value = 1275
$ws.Range("Q31").Value = 10
$ws.Range("R31").Value = "Hortaliza"
